# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that any leading "System"/"system" tokens are moved to the end of the
# comma-separated list, for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    # Collect the leading run of "System"/"system" tokens (case-insensitive).
    $lead = @()
    $i = 0
    while ($i -lt $parts.Length -and $parts[$i].ToLower() -eq "system") {
        $lead += $parts[$i]
        $i++
    }

    # Only reorder when there is a leading System-run AND other items remain.
    if ($lead.Length -gt 0 -and $i -lt $parts.Length) {
        $rest = $parts[$i..($parts.Length - 1)]
        $newParts = $rest + $lead
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value2 = $newVal
    }
}
